$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("size_correction")
$ws2 = $wb.Worksheets.Item("scale_decompression")

# New size-corrected delta values (AN:AS) for rows 3,4,5,6,7,8,9,10,11,12,13,14,15,16,17
$newData = @{
    3  = @(14.470670301233101, -2.7973453825958998, 17.268015683828999, 5.8366624593186298, 24.260364594243601, 47.550890704756597)
    4  = @(6.13555406578703, -12.722560186168201, 18.8581142519552, -3.2935030601905999, 17.6814351225267, 34.550395937069801)
    5  = @(-0.94059250953937301, 1.6218818131319399, -2.56247432267131, 0.340644651796284, 22.603707667668601, 44.269801263415197)
    6  = @(14.2561403318026, -1.0745010524815199, 15.3306413842841, 6.5908196396605696, 24.352520338578799, 47.733555896923399)
    7  = @(5.9213907120378897, -12.6401691821859, 18.561559894223802, -3.359389235074, 17.7807087203796, 34.745984622457499)
    8  = @(-1.58436029364894, -0.18319362963625699, -1.4011666640126801, -0.883776961642601, 21.5940953146835, 42.272658245997299)
    9  = @(15.3570170153258, -4.8295569856192602, 20.1865740009451, 5.2637300148532997, 23.307962421212899, 45.664003224382199)
    10 = @(13.189529284932201, -32.4605382585111, 45.650067543443299, -9.6355044867894595, 23.012027994814499, 45.078036862881099)
    11 = @(16.467420459830699, -4.6096001622838401, 21.0770206221145, 5.9289101487734399, 23.171290901832201, 45.393366379958998)
    12 = @(6.26118963247845, -13.3677308440912, 19.6289204765697, -3.5532706058064001, 17.047677096118001, 33.3021886649449)
    13 = @(-0.344641487913244, 0.61220577119724895, -0.95684725911049295, 0.133782141642002, 22.4590534325555, 43.983543024889798)
    14 = @(13.754238570479799, -3.6574638553349401, 17.411702425814799, 5.0483873575724498, 23.123577296198999, 45.298891917228502)
    15 = @(16.001773374926898, -4.3685511398915402, 20.3703245148184, 5.8166111175176898, 23.703515560060598, 46.447468335766501)
    16 = @(6.34330852432696, -13.6683362136055, 20.0116447379324, -3.66251384463928, 17.560418650252199, 34.311993678702798)
    17 = @(12.5964587639961, -26.429540998654598, 39.025999762650798, -6.9165411173292499, 22.706990474774301, 44.474212261421599)
}

$cols = @("AN", "AO", "AP", "AQ", "AR", "AS")
foreach ($row in $newData.Keys) {
    $vals = $newData[$row]
    for ($i = 0; $i -lt 6; $i++) {
        $ref = "$($cols[$i])$row"
        $ws1.Range($ref).Value = $vals[$i]
    }
}

# scale_decompression!B9:G11 now reference size_correction averages instead of hardcoded values
$ws2.Range("B9").Formula = "=AVERAGE(size_correction!AN3,size_correction!AN6,size_correction!AN9,size_correction!AN11,size_correction!AN14,size_correction!AN15)"
$ws2.Range("C9").Formula = "=AVERAGE(size_correction!AO3,size_correction!AO6,size_correction!AO9,size_correction!AO11,size_correction!AO14,size_correction!AO15)"
$ws2.Range("D9").Formula = "=AVERAGE(size_correction!AP3,size_correction!AP6,size_correction!AP9,size_correction!AP11,size_correction!AP14,size_correction!AP15)"
$ws2.Range("E9").Formula = "=AVERAGE(size_correction!AQ3,size_correction!AQ6,size_correction!AQ9,size_correction!AQ11,size_correction!AQ14,size_correction!AQ15)"
$ws2.Range("F9").Formula = "=AVERAGE(size_correction!AR3,size_correction!AR6,size_correction!AR9,size_correction!AR11,size_correction!AR14,size_correction!AR15)"
$ws2.Range("G9").Formula = "=AVERAGE(size_correction!AS3,size_correction!AS6,size_correction!AS9,size_correction!AS11,size_correction!AS14,size_correction!AS15)"

$ws2.Range("B10").Formula = "=AVERAGE(size_correction!AN5,size_correction!AN8,size_correction!AN13)"
$ws2.Range("C10").Formula = "=AVERAGE(size_correction!AO5,size_correction!AO8,size_correction!AO13)"
$ws2.Range("D10").Formula = "=AVERAGE(size_correction!AP5,size_correction!AP8,size_correction!AP13)"
$ws2.Range("E10").Formula = "=AVERAGE(size_correction!AQ5,size_correction!AQ8,size_correction!AQ13)"
$ws2.Range("F10").Formula = "=AVERAGE(size_correction!AR5,size_correction!AR8,size_correction!AR13)"
$ws2.Range("G10").Formula = "=AVERAGE(size_correction!AS5,size_correction!AS8,size_correction!AS13)"

$ws2.Range("B11").Formula = "=AVERAGE(size_correction!AN4,size_correction!AN7,size_correction!AN12,size_correction!AN16)"
$ws2.Range("C11").Formula = "=AVERAGE(size_correction!AO4,size_correction!AO7,size_correction!AO12,size_correction!AO16)"
$ws2.Range("D11").Formula = "=AVERAGE(size_correction!AP4,size_correction!AP7,size_correction!AP12,size_correction!AP16)"
$ws2.Range("E11").Formula = "=AVERAGE(size_correction!AQ4,size_correction!AQ7,size_correction!AQ12,size_correction!AQ16)"
$ws2.Range("F11").Formula = "=AVERAGE(size_correction!AR4,size_correction!AR7,size_correction!AR12,size_correction!AR16)"
$ws2.Range("G11").Formula = "=AVERAGE(size_correction!AS4,size_correction!AS7,size_correction!AS12,size_correction!AS16)"

# Update sheet view states
$ws1.Range("AP22").Select() | Out-Null
$ws2.Range("G25").Select() | Out-Null

$ws2.Activate()
